# Trade #116 closed at 2026-02-18 00:38:40 - unknown UNKNOWN +0.000%
# Plus two new OPEN trades (#173 momentum, #174 MarketMaking) logged afterwards.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet - aggregate totals move after trade #116 (HighProbConvergence)
#    closes and two new trades open.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.1     # Current Capital
$summary.Range("B4").Value = 0.21       # Total P&L $
$summary.Range("B6").Value = 144        # Total Trades
$summary.Range("B8").Value = 50         # Losing Trades
$summary.Range("B9").Value = 46.53      # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet - HighProbConvergence row (row 3) updates
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C3").Value = 100.39
$status.Range("D3").Value = 17
$status.Range("E3").Value = 0.39
$status.Range("F3").Value = 0.39
$status.Range("G3").Value = 64.70999999999999

# ---------------------------------------------------------------------------
# 3) All Trades sheet - trade #144 (row 145) closes
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(145, 7).Value = 0.081136          # Exit Price
$allTrades.Cells.Item(145, 8).Value = "CLOSED"          # Status
$allTrades.Cells.Item(145, 9).Value = -26.2401          # P&L %
$allTrades.Cells.Item(145, 10).Value = -0.03            # P&L $
$allTrades.Cells.Item(145, 11).Value = 100.39           # Capital After
$allTrades.Cells.Item(145, 12).Value = "early_exit"     # Exit Reason
$allTrades.Cells.Item(145, 13).Value = 0.1              # Duration (min)

# New trade rows appended to All Trades
$allTrades.Cells.Item(174, 1).Value = 173
$c = $allTrades.Cells.Item(174, 2)
$c.NumberFormat = "@"
$c.Value = "2026-02-18"
$c.Style = "Normal"
$allTrades.Cells.Item(174, 3).Value = "00:38:33"
$allTrades.Cells.Item(174, 4).Value = "momentum"
$allTrades.Cells.Item(174, 5).Value = "DOWN"
$allTrades.Cells.Item(174, 6).Value = 0.11
$allTrades.Cells.Item(174, 8).Value = "OPEN"
$allTrades.Cells.Item(174, 9).Value = 0
$allTrades.Cells.Item(174, 10).Value = 0
$allTrades.Cells.Item(174, 11).Value = 99.22374292899114
$allTrades.Cells.Item(174, 13).Value = 0
$allTrades.Cells.Item(174, 14).Value = 0
$allTrades.Cells.Item(174, 15).Value = 0
$allTrades.Cells.Item(174, 16).Value = 0.9
$allTrades.Cells.Item(174, 17).Value = "Downward momentum: -1.942% over 10 samples"

$allTrades.Cells.Item(175, 1).Value = 174
$c2 = $allTrades.Cells.Item(175, 2)
$c2.NumberFormat = "@"
$c2.Value = "2026-02-18"
$c2.Style = "Normal"
$allTrades.Cells.Item(175, 3).Value = "00:38:34"
$allTrades.Cells.Item(175, 4).Value = "MarketMaking"
$allTrades.Cells.Item(175, 5).Value = "DOWN"
$allTrades.Cells.Item(175, 6).Value = 0.1
$allTrades.Cells.Item(175, 8).Value = "OPEN"
$allTrades.Cells.Item(175, 9).Value = 0
$allTrades.Cells.Item(175, 10).Value = 0
$allTrades.Cells.Item(175, 11).Value = 99.21858346467945
$allTrades.Cells.Item(175, 13).Value = 0
$allTrades.Cells.Item(175, 14).Value = 0
$allTrades.Cells.Item(175, 15).Value = 0
$allTrades.Cells.Item(175, 16).Value = 0.6
$allTrades.Cells.Item(175, 17).Value = "Normal spread capture: 198 bps"

# ---------------------------------------------------------------------------
# 4) HighProbConvergence sheet - trade #144 (row 18) closes
# ---------------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Cells.Item(18, 7).Value = 0.081136
$hpc.Cells.Item(18, 8).Value = "CLOSED"
$hpc.Cells.Item(18, 9).Value = -26.2401
$hpc.Cells.Item(18, 10).Value = -0.03
$hpc.Cells.Item(18, 11).Value = 100.39
$hpc.Cells.Item(18, 16).Value = "early_exit"
$hpc.Cells.Item(18, 17).Value = 0.1

# ---------------------------------------------------------------------------
# 5) momentum sheet - new trade #173 (row 46)
# ---------------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Cells.Item(46, 1).Value = 173
$c3 = $momentum.Cells.Item(46, 2)
$c3.NumberFormat = "@"
$c3.Value = "2026-02-18"
$c3.Style = "Normal"
$momentum.Cells.Item(46, 3).Value = "00:38:33"
$momentum.Cells.Item(46, 4).Value = "momentum"
$momentum.Cells.Item(46, 5).Value = "DOWN"
$momentum.Cells.Item(46, 6).Value = 0.11
$momentum.Cells.Item(46, 8).Value = "OPEN"
$momentum.Cells.Item(46, 9).Value = 0
$momentum.Cells.Item(46, 10).Value = 0
$momentum.Cells.Item(46, 11).Value = 99.22374292899114
$momentum.Cells.Item(46, 12).Value = 0
$momentum.Cells.Item(46, 13).Value = 0
$momentum.Cells.Item(46, 14).Value = 0.9
$momentum.Cells.Item(46, 15).Value = "Downward momentum: -1.942% over 10 samples"
$momentum.Cells.Item(46, 17).Value = 0

# ---------------------------------------------------------------------------
# 6) MarketMaking sheet - new trade #174 (row 69)
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Cells.Item(69, 1).Value = 174
$c4 = $mm.Cells.Item(69, 2)
$c4.NumberFormat = "@"
$c4.Value = "2026-02-18"
$c4.Style = "Normal"
$mm.Cells.Item(69, 3).Value = "00:38:34"
$mm.Cells.Item(69, 4).Value = "MarketMaking"
$mm.Cells.Item(69, 5).Value = "DOWN"
$mm.Cells.Item(69, 6).Value = 0.1
$mm.Cells.Item(69, 8).Value = "OPEN"
$mm.Cells.Item(69, 9).Value = 0
$mm.Cells.Item(69, 10).Value = 0
$mm.Cells.Item(69, 11).Value = 99.21858346467945
$mm.Cells.Item(69, 12).Value = 0
$mm.Cells.Item(69, 13).Value = 0
$mm.Cells.Item(69, 14).Value = 0.6
$mm.Cells.Item(69, 15).Value = "Normal spread capture: 198 bps"
$mm.Cells.Item(69, 17).Value = 0
